$p = $ppt.ActivePresentation

# Slide 2, shape "Google Shape;219;p19" (10th shape) holds the title
# "Implementação da classe ListaCircular" split into two runs:
#   run 1: "Implementação da classe "
#   run 2: "ListaCircular" (with err="1" spell-flag)
# Target text: "Implementação da " + "classe No()"
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(10)

# First run: "Implementação da classe " (24 chars) -> "Implementação da "
$tr = $shp.TextFrame.TextRange
$run1 = $tr.Characters(1, 24)
$run1.Text = "Implementação da "

# Second run: now starts right after the (shorter) first run, 13 chars
# "ListaCircular" -> "classe No()"
$tr2 = $shp.TextFrame.TextRange
$run2 = $tr2.Characters(18, 13)
$run2.Text = "classe No()"
